$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.1 = 24670.06 pesos`n✅ 24670.06 pesos = 6.05 = 959.74 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 164
$wsTasas.Range("O10").Value = 4045.89
$wsTasas.Range("O12").Value = 158.53
